$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.3578751705562862, 3.024458662146624, 40.38354180535622, 6.354804623696642, 6.414829156315433, 46)
    3  = @(0.3027393802105286, 2.373015087887086, 23.34018936934246, 4.831168530422269, 4.876157719561091, 45)
    4  = @(0.1995050362595803, 2.425100243263678, 21.07208077481318, 4.590433615118857, 4.639116432198105, 44)
    5  = @(0.491114422551954, 2.30617629171989, 20.67429859926619, 4.546899888854623, 4.573795768127671, 43)
    6  = @(0.3200401550488365, 2.334673612171998, 21.38528161385673, 4.624422300553522, 4.669255825040207, 42)
    7  = @(0.5611040982166942, 2.339453670431203, 21.47116487205351, 4.633698832687933, 4.656740838783429, 41)
    8  = @(0.4062276606841116, 2.380866204400757, 22.22088215003843, 4.713903069648168, 4.756195553303624, 40)
    9  = @(0.5824843841155368, 2.368614004192091, 22.52413720137466, 4.745960092686691, 4.771651806576584, 39)
    10 = @(0.3786425555776947, 2.491417790537439, 23.33653245265391, 4.830790044356504, 4.880574196876067, 38)
    11 = @(0.6139585437634532, 2.407289227759336, 23.72539890019613, 4.870872498864667, 4.898675532022325, 37)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
